$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.621.71"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").Value = "3.059.58"
$ws.Range("E3").Value = "  +1.40%  "

# Row 5
$ws.Range("D5").Value = "514.74"
$ws.Range("E5").Value = "  +0.71%  "

# Row 6
$ws.Range("D6").Value = "140.38"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  -0.75%  "

# Row 9
$ws.Range("D9").Value = "7.21"
$ws.Range("E9").Value = "  -5.07%  "

# Row 10
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  -1.09%  "

# Row 11
$ws.Range("D11").Value = "0.377"
$ws.Range("E11").Value = "  +2.58%  "

# Row 12
$ws.Range("D12").Value = "3.584.82"
$ws.Range("E12").Value = "  +1.59%  "

# Row 13
$ws.Range("E13").Value = "  -3.42%  "

# Row 14
$ws.Range("D14").Value = "26.80"
$ws.Range("E14").Value = "  +0.34%  "

# Row 15
$ws.Range("E15").Value = "  +2.32%  "

# Row 16
$ws.Range("D16").Value = "57.596.59"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("D17").Value = "6.18"
$ws.Range("E17").Value = "  -0.93%  "

# Row 18
$ws.Range("D18").Value = "3.058.88"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +4.03%  "

# Row 20
$ws.Range("D20").Value = "8.17"
$ws.Range("E20").Value = "  +2.31%  "

# Row 21
$ws.Range("D21").Value = "330.71"
$ws.Range("E21").Value = "  -0.39%  "

# Row 22
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  +0.66%  "

# Row 23
$ws.Range("D23").Value = "0.507"
$ws.Range("E23").Value = "  +1.23%  "

# Row 24
$ws.Range("D24").Value = "65.68"
$ws.Range("E24").Value = "  +1.51%  "

# Row 25
$ws.Range("D25").Value = "3.186.35"
$ws.Range("E25").Value = "  +1.57%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  -3.16%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0909"
$ws.Range("E28").Value = "  -2.00%  "

# Row 29
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  -1.42%  "

# Row 30
$ws.Range("D30").Value = "7.24"
$ws.Range("E30").Value = "  -1.43%  "

# Row 31
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  +0.41%  "

# Row 33
$ws.Range("D33").Value = "20.84"
$ws.Range("E33").Value = "  +0.56%  "

# Row 34
$ws.Range("D34").Value = "154.08"
$ws.Range("E34").Value = "  -0.75%  "

# Row 35
$ws.Range("D35").Value = "4.67"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -0.80%  "

# Row 38
$ws.Range("D38").Value = "25.29"
$ws.Range("E38").Value = "  +2.84%  "

# Row 39
$ws.Range("D39").Value = "0.0680"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40
$ws.Range("D40").Value = "37.09"
$ws.Range("E40").Value = "  -1.43%  "

# Row 41
$ws.Range("E41").Value = "  +0.23%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.666"
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.40"
$ws.Range("E44").Value = "  -1.39%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.194.08"
$ws.Range("E45").Value = "  -1.89%  "

# Row 46
$ws.Range("D46").Value = "6.10"
$ws.Range("E46").Value = "  +1.10%  "

# Row 47
$ws.Range("D47").Value = "0.954"
$ws.Range("E47").Value = "  -3.63%  "

# Row 48
$ws.Range("E48").Value = "  +1.99%  "

# Row 49
$ws.Range("D49").Value = "20.08"
$ws.Range("E49").Value = "  +2.69%  "

# Row 50
$ws.Range("B50").Value = "Notcoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/2L2Y4ghjj+notcoin-not"
$ws.Range("D50").Value = "0.0172"
$ws.Range("E50").Value = "  +8.52%  "

# Row 51
$ws.Range("D51").Value = "0.184"
$ws.Range("E51").Value = "  -0.51%  "
